$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Restore value for rule R30 "From" (C10) from 18 to 1
$ws.Range("C10").Value = 1
